$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the activity log with a new row (row 7), reusing the formatting
# of the last existing data row (row 6) so the new date/time cells pick up
# the same date & time-of-day number formats already used in the table.
$ws.Rows("6").Copy()
$ws.Rows("7").Insert(-4121)   # xlShiftDown

# Fill in the new row's values
$ws.Range("A7").Value = "Chapter 1 still on process"
$ws.Range("B7").Value = 43755
$ws.Range("C7").Value = 43755
$ws.Range("D7").Value = 0.83333333333333337
$ws.Range("E7").Value = 0.54166666666666663

# Update the view: zoom level and current selection moved to D8
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("D8").Select()
